# Remove the first course entry (and its blank separator row) from each of
# the 7 course worksheets. Each sheet is laid out as:
#   row 1      : header ("Course Title" / "Course Description")
#   row 2      : blank separator
#   row 3      : first course
#   row 4      : blank separator
#   row 5      : second course
#   ...
# Deleting rows 2:3 removes the first course + its separator and shifts
# everything else up by two rows, shrinking the used range by 2 rows.

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "Systems Management & Security",
    "Web Programming",
    "Programming",
    "IT Generalist",
    "Embedded Systems",
    "Database Admin",
    "Data Analytics"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows("2:3").Delete()
}

# The hidden ExternalData_N defined names (one per sheet, pointing at the
# query table's result range) track the query refresh range and are not
# automatically resized by the row delete above, so update them explicitly
# to match the new (shorter) used range of each sheet.
$wb.Names.Item("Data Analytics!ExternalData_1").RefersTo = "='Data Analytics'!`$A`$1:`$B`$61"
$wb.Names.Item("Database Admin!ExternalData_2").RefersTo = "='Database Admin'!`$A`$1:`$B`$61"
$wb.Names.Item("Embedded Systems!ExternalData_3").RefersTo = "='Embedded Systems'!`$A`$1:`$B`$31"
$wb.Names.Item("IT Generalist!ExternalData_4").RefersTo = "='IT Generalist'!`$A`$1:`$B`$57"
$wb.Names.Item("Programming!ExternalData_5").RefersTo = "=Programming!`$A`$1:`$B`$51"
$wb.Names.Item("Web Programming!ExternalData_6").RefersTo = "='Web Programming'!`$A`$1:`$B`$51"
$wb.Names.Item("Systems Management & Security!ExternalData_7").RefersTo = "='Systems Management & Security'!`$A`$1:`$B`$59"

# Restore cell selections to match the post-edit active cells.
$ws1 = $wb.Worksheets.Item("Systems Management & Security")
$ws1.Activate()
$ws1.Range("D4").Select()

$ws2 = $wb.Worksheets.Item("Web Programming")
$ws2.Activate()
$ws2.Range("B2").Select()

$ws3 = $wb.Worksheets.Item("Programming")
$ws3.Activate()
$ws3.Range("B2").Select()

$ws4 = $wb.Worksheets.Item("IT Generalist")
$ws4.Activate()
$ws4.Range("A2:A3").Select()

$ws5 = $wb.Worksheets.Item("Embedded Systems")
$ws5.Activate()
$ws5.Range("A4").Select()

$ws6 = $wb.Worksheets.Item("Database Admin")
$ws6.Activate()
$ws6.Range("A2:A3").Select()

$ws7 = $wb.Worksheets.Item("Data Analytics")
$ws7.Activate()
$ws7.Range("E4").Select()
